# update electrolysis data for methanol case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# Row 2 (Solar_Plant_Kasso): drop the ramp/start-up/shut-down figures that no
# longer apply, keep Cap_Output1_existing (AB2) as-is
$ws.Range("O2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3 (Electrolyzer): refresh the electrolysis figures
$ws.Range("W3").Value = 0.0063
$ws.Range("Y3").Value = 1.76
$ws.Range("AD3").ClearContents()
$ws.Range("AB3").Value = 216.9

# Bring the sheet to the front and move the selection to where the author
# left it
$ws.Activate()
$ws.Range("X24").Select()
